$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve Price column values as exact text (avoid numeric/date coercion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.412.36'
$ws.Range("E2").Value = '  +0.83%  '

$ws.Range("D3").Value = '1.867.81'
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("D4").Value = '1.016'
$ws.Range("E4").Value = '  +1.35%  '

$ws.Range("D5").Value = '316.42'
$ws.Range("E5").Value = '  +1.01%  '

$ws.Range("D6").Value = '1.017'
$ws.Range("E6").Value = '  +1.49%  '

$ws.Range("D7").Value = '0.5131'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '0.3926'
$ws.Range("E8").Value = '  +0.59%  '

$ws.Range("D9").Value = '0.08303'
$ws.Range("E9").Value = '  -0.87%  '

$ws.Range("D10").Value = '1.116'
$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("D11").Value = '41.98'
$ws.Range("E11").Value = '  +1.26%  '

$ws.Range("D12").Value = '6.236'
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").Value = '1.864.74'
$ws.Range("E13").Value = '  -0.58%  '

$ws.Range("D14").Value = '20.28'
$ws.Range("E14").Value = '  -2.22%  '

$ws.Range("D15").Value = '7.205'
$ws.Range("E15").Value = '  -1.27%  '

$ws.Range("D16").Value = '1.016'
$ws.Range("E16").Value = '  +1.34%  '

$ws.Range("D17").Value = '0.00001102'
$ws.Range("E17").Value = '  -0.46%  '

$ws.Range("D18").Value = '91.23'
$ws.Range("E18").Value = '  -0.01%  '

$ws.Range("D19").Value = '0.06732'
$ws.Range("E19").Value = '  +1.12%  '

$ws.Range("D20").Value = '17.63'
$ws.Range("E20").Value = '  -0.58%  '

$ws.Range("D21").Value = '1.017'
$ws.Range("E21").Value = '  +1.54%  '

$ws.Range("D22").Value = '5.954'
$ws.Range("E22").Value = '  -1.62%  '

$ws.Range("D23").Value = '28.450.99'
$ws.Range("E23").Value = '  +0.83%  '

$ws.Range("D24").Value = '11.11'
$ws.Range("E24").Value = '  -0.75%  '

$ws.Range("D25").Value = '2.277'
$ws.Range("E25").Value = '  +0.69%  '

$ws.Range("D26").Value = '2.078.63'
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("D27").Value = '161.51'
$ws.Range("E27").Value = '  +1.96%  '

$ws.Range("D28").Value = '20.68'
$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").Value = '2.407'
$ws.Range("E29").Value = '  -3.79%  '

$ws.Range("D30").Value = '126.55'
$ws.Range("E30").Value = '  +1.04%  '

$ws.Range("D31").Value = '0.1042'
$ws.Range("E31").Value = '  -2.15%  '

$ws.Range("D32").Value = '1.033'
$ws.Range("E32").Value = '  -0.78%  '

$ws.Range("D33").Value = '5.809'
$ws.Range("E33").Value = '  -1.50%  '

$ws.Range("D34").Value = '3.637'
$ws.Range("E34").Value = '  +1.31%  '

$ws.Range("D35").Value = '0.02440'
$ws.Range("E35").Value = '  -0.51%  '

$ws.Range("D36").Value = '0.06485'
$ws.Range("E36").Value = '  -1.05%  '

$ws.Range("D37").Value = '9.125'
$ws.Range("E37").Value = '  -6.35%  '

$ws.Range("D38").Value = '0.2175'
$ws.Range("E38").Value = '  -0.94%  '

$ws.Range("D39").Value = '1.260'
$ws.Range("E39").Value = '  +2.35%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '1.183'
$ws.Range("E40").Value = '  -2.20%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6423'
$ws.Range("E41").Value = '  -1.61%  '

$ws.Range("D42").Value = '4.966'
$ws.Range("E42").Value = '  -1.12%  '

$ws.Range("D43").Value = '11.10'
$ws.Range("E43").Value = '  -1.92%  '

$ws.Range("D44").Value = '0.5994'
$ws.Range("E44").Value = '  -2.12%  '

$ws.Range("D45").Value = '12.99'

$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.292'
$ws.Range("E46").Value = '  +0.27%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.707'
$ws.Range("E47").Value = '  +1.01%  '

$ws.Range("D48").Value = '1.986'
$ws.Range("E48").Value = '  -1.48%  '

$ws.Range("D49").Value = '1.201'
$ws.Range("E49").Value = '  -1.95%  '

$ws.Range("D50").Value = '121.06'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("D51").Value = '0.06863'
$ws.Range("E51").Value = '  -0.56%  '
